$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.088.38"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "3.100.28"
$ws.Range("E3").Value = "  +4.99%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'583.27"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Value = "'169.28"
$ws.Range("E6").Value = "  +7.12%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.094.65"
$ws.Range("E8").Value = "  +4.86%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").Value = "'6.68"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D12").Value = "'0.482"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "'36.91"
$ws.Range("E14").Value = "  +8.54%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.616.07"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "67.122.64"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D19").Value = "3.103.35"
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("D20").Value = "'16.28"
$ws.Range("E20").Value = "  +17.52%  "
$ws.Range("D21").Value = "'471.39"
$ws.Range("E21").Value = "  +5.41%  "
$ws.Range("D22").Value = "'0.716"
$ws.Range("E22").Value = "  +5.39%  "
$ws.Range("D23").Value = "'7.55"
$ws.Range("E23").Value = "  +4.59%  "
$ws.Range("D24").Value = "'83.91"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +8.80%  "
$ws.Range("D26").Value = "'12.94"
$ws.Range("E26").Value = "  +7.47%  "
$ws.Range("D27").Value = "'10.22"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'8.08"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("E30").Value = "  +5.26%  "
$ws.Range("D31").Value = "'2.68"
$ws.Range("E31").Value = "  +4.67%  "
$ws.Range("E32").Value = "  +4.51%  "
$ws.Range("D33").Value = "'28.46"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("D38").Value = "'47.26"
$ws.Range("E38").Value = "  +10.20%  "
$ws.Range("D39").Value = "'2.10"
$ws.Range("E39").Value = "  +6.54%  "
$ws.Range("E40").Value = "  +7.09%  "
$ws.Range("D41").Value = "'50.41"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.124"
$ws.Range("E42").Value = "  +4.40%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.92"
$ws.Range("E43").Value = "  +5.39%  "
$ws.Range("D44").Value = "'8.77"
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("D45").Value = "'395.12"
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("D46").Value = "'0.0364"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").Value = "2.766.50"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "'135.08"
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'24.87"
$ws.Range("E50").Value = "  +7.29%  "
$ws.Range("E51").Value = "  +4.81%  "
